# Update existing rows 2-10 and append new rows 11-16 with refreshed NATMI TPM output
# (adds Inflammatory-Mac / Resolving-Mac target clusters and recomputed scores)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt9a"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6558959999999999
$ws.Range("H2").Value = 1.967688
$ws.Range("I2").Value = 0.1303533075142512
$ws.Range("J2").Value = 0.1303533075142512
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 29.75868033333333
$ws.Range("N2").Value = 89.27604099999999
$ws.Range("O2").Value = 0.4948552779010537
$ws.Range("P2").Value = 0.4948552779010535
$ws.Range("Q2").Value = 19.51859939591199
$ws.Range("R2").Value = 175.667394563208
$ws.Range("S2").Value = 0.06450602221528628
$ws.Range("T2").Value = 0.06450602221528627

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt9a"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6558959999999999
$ws.Range("H3").Value = 1.967688
$ws.Range("I3").Value = 0.1303533075142512
$ws.Range("J3").Value = 0.1303533075142512
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.55525033333333
$ws.Range("N3").Value = 52.665751
$ws.Range("O3").Value = 0.2919251856942525
$ws.Range("P3").Value = 0.2919251856942524
$ws.Range("Q3").Value = 11.514418472632
$ws.Range("R3").Value = 103.629766253688
$ws.Range("S3").Value = 0.03805341350195778
$ws.Range("T3").Value = 0.03805341350195777

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt9a"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6558959999999999
$ws.Range("H4").Value = 1.967688
$ws.Range("I4").Value = 0.1303533075142512
$ws.Range("J4").Value = 0.1303533075142512
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1421396666666667
$ws.Range("N4").Value = 0.426419
$ws.Range("O4").Value = 0.002363631836533717
$ws.Range("P4").Value = 0.002363631836533717
$ws.Range("Q4").Value = 0.09322883880799998
$ws.Range("R4").Value = 0.8390595492719999
$ws.Range("S4").Value = 0.0003081072276381539
$ws.Range("T4").Value = 0.0003081072276381539

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt9a"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6558959999999999
$ws.Range("H5").Value = 1.967688
$ws.Range("I5").Value = 0.1303533075142512
$ws.Range("J5").Value = 0.1303533075142512
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.42872866666667
$ws.Range("N5").Value = 37.286186
$ws.Range("O5").Value = 0.2066765699758167
$ws.Range("P5").Value = 0.2066765699758166
$ws.Range("Q5").Value = 8.151953417551999
$ws.Range("R5").Value = 73.367580757968
$ws.Range("S5").Value = 0.02694097448204829
$ws.Range("T5").Value = 0.02694097448204828

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt9a"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6558959999999999
$ws.Range("H6").Value = 1.967688
$ws.Range("I6").Value = 0.1303533075142512
$ws.Range("J6").Value = 0.1303533075142512
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.251329
$ws.Range("N6").Value = 0.753987
$ws.Range("O6").Value = 0.004179334592343558
$ws.Range("P6").Value = 0.004179334592343557
$ws.Range("Q6").Value = 0.164845685784
$ws.Range("R6").Value = 1.483611172056
$ws.Range("S6").Value = 0.0005447900873207075
$ws.Range("T6").Value = 0.0005447900873207074

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt9a"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.335609
$ws.Range("H7").Value = 10.006827
$ws.Range("I7").Value = 0.6629216609406124
$ws.Range("J7").Value = 0.6629216609406123
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 29.75868033333333
$ws.Range("N7").Value = 89.27604099999999
$ws.Range("O7").Value = 0.4948552779010537
$ws.Range("P7").Value = 0.4948552779010535
$ws.Range("Q7").Value = 99.26332194798965
$ws.Range("R7").Value = 893.3698975319069
$ws.Range("S7").Value = 0.3280502827513949
$ws.Range("T7").Value = 0.3280502827513947

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt9a"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.335609
$ws.Range("H8").Value = 10.006827
$ws.Range("I8").Value = 0.6629216609406124
$ws.Range("J8").Value = 0.6629216609406123
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 17.55525033333333
$ws.Range("N8").Value = 52.665751
$ws.Range("O8").Value = 0.2919251856942525
$ws.Range("P8").Value = 0.2919251856942524
$ws.Range("Q8").Value = 58.55745100911967
$ws.Range("R8").Value = 527.017059082077
$ws.Range("S8").Value = 0.1935235289708306
$ws.Range("T8").Value = 0.1935235289708305

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt9a"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.335609
$ws.Range("H9").Value = 10.006827
$ws.Range("I9").Value = 0.6629216609406124
$ws.Range("J9").Value = 0.6629216609406123
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1421396666666667
$ws.Range("N9").Value = 0.426419
$ws.Range("O9").Value = 0.002363631836533717
$ws.Range("P9").Value = 0.002363631836533717
$ws.Range("Q9").Value = 0.4741223513903333
$ws.Range("R9").Value = 4.267101162513
$ws.Range("S9").Value = 0.001566902742927042
$ws.Range("T9").Value = 0.001566902742927041

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt9a"
$ws.Range("C10").Value = "Fzd4"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.335609
$ws.Range("H10").Value = 10.006827
$ws.Range("I10").Value = 0.6629216609406124
$ws.Range("J10").Value = 0.6629216609406123
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.42872866666667
$ws.Range("N10").Value = 37.286186
$ws.Range("O10").Value = 0.2066765699758167
$ws.Range("P10").Value = 0.2066765699758166
$ws.Range("Q10").Value = 41.45737919909133
$ws.Range("R10").Value = 373.116412791822
$ws.Range("S10").Value = 0.1370103750458771
$ws.Range("T10").Value = 0.137010375045877

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt9a"
$ws.Range("C11").Value = "Fzd4"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.335609
$ws.Range("H11").Value = 10.006827
$ws.Range("I11").Value = 0.6629216609406124
$ws.Range("J11").Value = 0.6629216609406123
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.251329
$ws.Range("N11").Value = 0.753987
$ws.Range("O11").Value = 0.004179334592343558
$ws.Range("P11").Value = 0.004179334592343557
$ws.Range("Q11").Value = 0.8383352743609999
$ws.Range("R11").Value = 7.545017469248999
$ws.Range("S11").Value = 0.002770571429582949
$ws.Range("T11").Value = 0.002770571429582948

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Wnt9a"
$ws.Range("C12").Value = "Fzd4"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.040174
$ws.Range("H12").Value = 3.120522
$ws.Range("I12").Value = 0.2067250315451364
$ws.Range("J12").Value = 0.2067250315451363
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 29.75868033333333
$ws.Range("N12").Value = 89.27604099999999
$ws.Range("O12").Value = 0.4948552779010537
$ws.Range("P12").Value = 0.4948552779010535
$ws.Range("Q12").Value = 30.95420555704467
$ws.Range("R12").Value = 278.587850013402
$ws.Range("S12").Value = 0.1022989729343726
$ws.Range("T12").Value = 0.1022989729343725

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Wnt9a"
$ws.Range("C13").Value = "Fzd4"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.040174
$ws.Range("H13").Value = 3.120522
$ws.Range("I13").Value = 0.2067250315451364
$ws.Range("J13").Value = 0.2067250315451363
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 17.55525033333333
$ws.Range("N13").Value = 52.665751
$ws.Range("O13").Value = 0.2919251856942525
$ws.Range("P13").Value = 0.2919251856942524
$ws.Range("Q13").Value = 18.26051496022467
$ws.Range("R13").Value = 164.344634642022
$ws.Range("S13").Value = 0.06034824322146413
$ws.Range("T13").Value = 0.06034824322146411

# Row 14
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Wnt9a"
$ws.Range("C14").Value = "Fzd4"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.040174
$ws.Range("H14").Value = 3.120522
$ws.Range("I14").Value = 0.2067250315451364
$ws.Range("J14").Value = 0.2067250315451363
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1421396666666667
$ws.Range("N14").Value = 0.426419
$ws.Range("O14").Value = 0.002363631836533717
$ws.Range("P14").Value = 0.002363631836533717
$ws.Range("Q14").Value = 0.1478499856353334
$ws.Range("R14").Value = 1.330649870718
$ws.Range("S14").Value = 0.0004886218659685213
$ws.Range("T14").Value = 0.000488621865968521

# Row 15
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Wnt9a"
$ws.Range("C15").Value = "Fzd4"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.040174
$ws.Range("H15").Value = 3.120522
$ws.Range("I15").Value = 0.2067250315451364
$ws.Range("J15").Value = 0.2067250315451363
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 12.42872866666667
$ws.Range("N15").Value = 37.286186
$ws.Range("O15").Value = 0.2066765699758167
$ws.Range("P15").Value = 0.2066765699758166
$ws.Range("Q15").Value = 12.92804041212134
$ws.Range("R15").Value = 116.352363709092
$ws.Range("S15").Value = 0.04272522044789129
$ws.Range("T15").Value = 0.04272522044789127

# Row 16
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Wnt9a"
$ws.Range("C16").Value = "Fzd4"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.040174
$ws.Range("H16").Value = 3.120522
$ws.Range("I16").Value = 0.2067250315451364
$ws.Range("J16").Value = 0.2067250315451363
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.251329
$ws.Range("N16").Value = 0.753987
$ws.Range("O16").Value = 0.004179334592343558
$ws.Range("P16").Value = 0.004179334592343557
$ws.Range("Q16").Value = 0.261425891246
$ws.Range("R16").Value = 2.352833021214
$ws.Range("S16").Value = 0.0008639730754399017
$ws.Range("T16").Value = 0.0008639730754399013

